$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player added to Adelaide Strikers (column A), row 21
$ws.Range("A21").Value = "Tom Andrews"

# New player added to Sydney Thunder (column H), row 20
$ws.Range("H20").Value = "Gurinder Sandu"

# Update view: scroll back to top-left A1, select H20
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("H20").Select()

# Page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
